$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.843326807022095
$ws.Range("B1").Value = 2.897255182266235
$ws.Range("C1").Value = 1.928537368774414
$ws.Range("D1").Value = 1.690407752990723
$ws.Range("E1").Value = 1.671060681343079
